$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for the A (Id), I (Antal), Q (Ost) and R (Nord)
# columns of rows 35-47 before any changes are made, so the permutation
# below can be applied consistently using a single source snapshot.
$rows = 35..47
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        A = $ws.Range("A$r").Value2
        I = $ws.Range("I$r").Value2
        Q = $ws.Range("Q$r").Value2
        R = $ws.Range("R$r").Value2
    }
}

# Mapping of destination row -> source row (values copied from the
# source row's original data into the destination row).
$mapping = @{
    35 = 38
    36 = 37
    37 = 41
    38 = 42
    39 = 40
    40 = 39
    41 = 36
    42 = 44
    43 = 46
    44 = 35
    46 = 47
    47 = 43
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]

    foreach ($col in @("A", "I", "Q", "R")) {
        $newVal = $orig[$src][$col]
        $curVal = $orig[$dest][$col]
        # Only touch the cell when the value actually changes; this avoids
        # rewriting (and thereby silently re-typing) cells whose value is
        # coincidentally identical between source and destination rows
        # (e.g. I41/I44, which both stay "1").
        if ($curVal -ne $newVal) {
            $ws.Range("$col$dest").Value2 = $newVal
        }
    }
}
